$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (dbExcel / Neo4j filename) rows 2-4 ---
$origD = @'
TC01_CDS_Filter_InstrumentModel-DNBSEQ-G400_Neo4jData.xlsx
'@
$newD = $origD.Replace("TC01", "TC05").Replace("DNBSEQ-G400", "Illumina HiSeq 2000")
$ws.Range("D2:D4").Value = $newD

# --- Column E (WebExcel filename) rows 2-4 ---
$origE = @'
TC01_CDS_Filter_InstrumentModel-DNBSEQ-G400_WebData.xlsx
'@
$newE = $origE.Replace("TC01", "TC05").Replace("DNBSEQ-G400", "Illumina HiSeq 2000")
$ws.Range("E2:E4").Value = $newE

# --- Column B (query text), each row has a different query ---
$origB2 = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['DNBSEQ-G400']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY `Participant ID`LIMIT 100
'@
$ws.Range("B2").Value = $origB2.Replace("DNBSEQ-G400", "Illumina HiSeq 2000")

$origB3 = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['DNBSEQ-G400']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@
$ws.Range("B3").Value = $origB3.Replace("DNBSEQ-G400", "Illumina HiSeq 2000")

$origB4 = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['DNBSEQ-G400']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@
$ws.Range("B4").Value = $origB4.Replace("DNBSEQ-G400", "Illumina HiSeq 2000")

# --- Column C (StatQuery), same text for rows 2-4 ---
$origC = @'
MATCH (f:file)
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['DNBSEQ-G400']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,f, s, collect(distinct samp.sample_id) as samp
RETURN
count(distinct s) AS Studies,
count(distinct p) AS Participants,
count(distinct samp) AS Samples,
count(distinct f) AS Files
'@
$newC = $origC.Replace("DNBSEQ-G400", "Illumina HiSeq 2000")
$ws.Range("C2:C4").Value = $newC

# --- Column widths (bestFit autofit changed due to new content) ---
$ws.Columns.Item(4).ColumnWidth = 94.85546875
$ws.Columns.Item(5).ColumnWidth = 93.140625

# --- Selection change ---
$ws.Range("D3").Select()
